$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to Text format so Excel
# keeps them as strings (matching the source inline-string cells) instead of
# auto-converting them to numbers.

$ws.Range("D2").Value = "25.620.62"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "1.670.27"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.80"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4813"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2638"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06170"
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07093"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "1.664.89"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.92"
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6011"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.426"
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.82"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9985"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9987"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "25.606.86"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006804"
$ws.Range("E19").Value = "  +4.46%  "
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.496"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "1.877.58"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.740"
$ws.Range("E23").Value = "  +3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.392"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.68"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.407"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.87"
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.711"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.983"
$ws.Range("E30").Value = "  +4.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.687"
$ws.Range("E31").Value = "  +4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07693"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04387"
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9975"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.619"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6211"
$ws.Range("E36").Value = "  +6.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9551"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.625"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8706"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9983"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01518"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.24"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3797"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.678"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1126"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.261"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05261"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.68"
$ws.Range("E49").Value = "  +0.84%  "

# Row 50 / 51: Decentraland and EnergySwap swap list position, each taking new
# Price / Volume(1h) values (not a simple swap of the old values).
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3363"
$ws.Range("E50").Value = "  +2.00%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.355"
$ws.Range("E51").Value = "  +0.35%  "
